# Version archivos 10 junio 2020
# Mark three additional students with an "X" in the "Asistencia" sheet's
# date column (column D). Two of them (rows 9 and 14) are flagged in red
# to call out a new/late addition, the third (row 37) uses the sheet's
# normal centered style, matching the rest of the "X" marks in the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asistencia")

# Row 9 -> D9: new mark, red font
$ws.Range("D9").Value = "X"
$ws.Range("D9").Font.Color = 255

# Row 14 -> D14: new mark, red font
$ws.Range("D14").Value = "X"
$ws.Range("D14").Font.Color = 255

# Row 37 -> D37: new mark, default (black) style like the other entries
$ws.Range("D37").Value = "X"

# Leave the selection on the last edited cell, like the saved workbook.
$ws.Range("D37").Select()
